$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking row)
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total row)
$ws.Range("B12").Value = 120
$ws.Range("C12").Value = -2.4
$ws.Range("E12").Value = "117.6/140"
